# Append new EURUSD daily/intraday OHLCV rows (4993-5057) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new rows with the same row style as the last existing row (row 4992),
# then overwrite with the actual imported values.
$ws.Range("A4992:F4992").Copy($ws.Range("A4993:F5057"))

$newData = @(
    @(45229, 1.05626, 1.05707, 1.05551, 1.05625, 11757),
    @(45229, 1.05625, 1.05658, 1.05583, 1.05609, 15372),
    @(45229, 1.05609, 1.05674, 1.05471, 1.05598, 29587),
    @(45229, 1.05598, 1.06084, 1.05597, 1.06068, 41413),
    @(45230, 1.06068, 1.06253, 1.05908, 1.06071, 74171),
    @(45230, 1.06071, 1.06217, 1.06001, 1.0615, 22570),
    @(45230, 1.0615, 1.06172, 1.06093, 1.06115, 6936),
    @(45230, 1.06115, 1.06126, 1.05906, 1.05938, 28030),
    @(45230, 1.05938, 1.06458, 1.05926, 1.06453, 34443),
    @(45230, 1.06453, 1.06749, 1.06211, 1.06241, 60573),
    @(45231, 1.06241, 1.06257, 1.05574, 1.05675, 90552),
    @(45231, 1.05675, 1.05858, 1.05657, 1.05744, 29256),
    @(45231, 1.05744, 1.05806, 1.05685, 1.05733, 10055),
    @(45231, 1.05733, 1.05764, 1.05658, 1.05672, 18408),
    @(45231, 1.05672, 1.05775, 1.0557, 1.05582, 28917),
    @(45231, 1.05582, 1.05627, 1.05369, 1.05385, 47594),
    @(45232, 1.05385, 1.05708, 1.05215, 1.05338, 80800),
    @(45232, 1.05338, 1.05758, 1.05165, 1.05698, 91229),
    @(45232, 1.05698, 1.0602, 1.05645, 1.06017, 14401),
    @(45232, 1.06017, 1.06021, 1.05932, 1.05959, 17839),
    @(45232, 1.05959, 1.06103, 1.05911, 1.06005, 33717),
    @(45232, 1.06005, 1.06675, 1.0599, 1.06641, 59586),
    @(45233, 1.06641, 1.06668, 1.06104, 1.06153, 74056),
    @(45233, 1.06153, 1.06325, 1.06137, 1.06216, 27223),
    @(45233, 1.06216, 1.06237, 1.06145, 1.06181, 9552),
    @(45233, 1.06181, 1.06287, 1.06147, 1.06263, 14019),
    @(45233, 1.06263, 1.0634, 1.06197, 1.0627, 23544),
    @(45233, 1.0627, 1.07194, 1.06241, 1.07093, 59851),
    @(45234, 1.07093, 1.07393, 1.06935, 1.07392, 106662),
    @(45234, 1.07392, 1.07467, 1.07224, 1.07277, 26350),
    @(45236, 1.07215, 1.07335, 1.07215, 1.07283, 13826),
    @(45236, 1.07283, 1.07388, 1.07263, 1.07355, 14112),
    @(45236, 1.07355, 1.07563, 1.07315, 1.07471, 35364),
    @(45236, 1.07471, 1.07557, 1.0735, 1.07374, 34867),
    @(45237, 1.07374, 1.07493, 1.0726, 1.07354, 43962),
    @(45237, 1.07354, 1.0737, 1.07163, 1.07172, 22305),
    @(45237, 1.07172, 1.07221, 1.07079, 1.07103, 12130),
    @(45237, 1.07103, 1.0714, 1.07052, 1.07104, 16241),
    @(45237, 1.07104, 1.07116, 1.0683, 1.06836, 34146),
    @(45237, 1.06836, 1.06976, 1.06655, 1.06759, 41765),
    @(45238, 1.06759, 1.06928, 1.0664, 1.06773, 52990),
    @(45238, 1.06773, 1.07048, 1.0676, 1.06997, 24363),
    @(45238, 1.06997, 1.06998, 1.06863, 1.06918, 13102),
    @(45238, 1.06918, 1.06991, 1.06853, 1.06861, 15441),
    @(45238, 1.06861, 1.0688, 1.06602, 1.06651, 33518),
    @(45238, 1.06651, 1.06778, 1.06592, 1.06692, 35547),
    @(45239, 1.06692, 1.0716, 1.06689, 1.07077, 55068),
    @(45239, 1.07077, 1.07141, 1.06999, 1.0709, 26556),
    @(45239, 1.0709, 1.07146, 1.07023, 1.0712, 12311),
    @(45239, 1.0712, 1.07156, 1.07021, 1.07025, 13604),
    @(45239, 1.07025, 1.07097, 1.06934, 1.06977, 31167),
    @(45239, 1.06977, 1.07088, 1.068, 1.07075, 39342),
    @(45240, 1.07075, 1.07255, 1.06998, 1.07047, 56595),
    @(45240, 1.07047, 1.07056, 1.06601, 1.06673, 64566),
    @(45240, 1.06673, 1.06704, 1.06641, 1.06653, 3460),
    @(45240, 1.07075, 1.07255, 1.06998, 1.07047, 56603),
    @(45240, 1.07047, 1.07056, 1.06601, 1.06673, 64726),
    @(45240, 1.06673, 1.06751, 1.06629, 1.0672, 13983),
    @(45240, 1.0672, 1.06744, 1.06657, 1.0669, 13820),
    @(45240, 1.0669, 1.068, 1.06561, 1.06772, 35112),
    @(45240, 1.06772, 1.06928, 1.06693, 1.06887, 38382),
    @(45241, 1.06887, 1.06896, 1.06623, 1.06792, 56460),
    @(45241, 1.06792, 1.06875, 1.06729, 1.06842, 20019),
    @(45243, 1.0683, 1.06907, 1.06805, 1.06867, 13886),
    @(45243, 1.06867, 1.06873, 1.06814, 1.06855, 6473)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = 4993 + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
